$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4038.6
$ws.Range("I17").Value = 3000.25
$ws.Range("J17").Value = 4730.8335
$ws.Range("K17").Value = 9000.75
$ws.Range("L17").Value = 14192.5005
$ws.Range("M17").Value = -8832.75
$ws.Range("N17").Value = -14528.5005

$ws.Range("H112").Value = 2120.9473
$ws.Range("I112").Value = 1499.5
$ws.Range("K112").Value = 4498.5
$ws.Range("M112").Value = -3390.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6362.6875
$ws.Range("I32").Value = 4955.032
$ws.Range("K32").Value = 4955.032
$ws.Range("M32").Value = -4668.032

$ws.Range("H107").Value = 20228
$ws.Range("J107").Value = 20228
$ws.Range("L107").Value = 20228
$ws.Range("N107").Value = -27908

$ws.Range("H132").Value = 2479.6667
$ws.Range("I132").Value = 2616.7144
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 7850.1432
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -5320.1432
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2082.25
$ws.Range("I86").Value = 2255.4285
$ws.Range("J86").Value = 1839.8
$ws.Range("K86").Value = 2255.4285
$ws.Range("L86").Value = 1839.8
$ws.Range("M86").Value = -1132.4285
$ws.Range("N86").Value = -4085.8

$ws.Range("H89").Value = 2082.25
$ws.Range("I89").Value = 2255.4285
$ws.Range("J89").Value = 1839.8
$ws.Range("K89").Value = 11277.1425
$ws.Range("L89").Value = 9199
$ws.Range("M89").Value = -5661.1425
$ws.Range("N89").Value = -20431

$ws.Range("H92").Value = 34999
$ws.Range("J92").Value = 34999
$ws.Range("L92").Value = 34999
$ws.Range("N92").Value = -39991

$ws.Range("H107").Value = 1758.8334
$ws.Range("I107").Value = 1110.6
$ws.Range("K107").Value = 1110.6
$ws.Range("M107").Value = 809.4000000000001

$ws.Range("H134").Value = 13940
$ws.Range("J134").Value = 10000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1860.375
$ws.Range("I3").Value = 418
$ws.Range("J3").Value = 3302.75
$ws.Range("K3").Value = 418
$ws.Range("L3").Value = 3302.75
$ws.Range("M3").Value = -305
$ws.Range("N3").Value = -3528.75

$ws.Range("H36").Value = 19444
$ws.Range("I36").Value = 19444
$ws.Range("K36").Value = 19444
$ws.Range("M36").Value = -19056

$ws.Range("H40").Value = 19444
$ws.Range("I40").Value = 19444
$ws.Range("K40").Value = 19444
$ws.Range("M40").Value = -19284

$ws.Range("H103").Value = 7074.5
$ws.Range("I103").Value = 7074.5
$ws.Range("K103").Value = 7074.5
$ws.Range("M103").Value = -5902.5

$ws.Range("H132").Value = 1996.5
$ws.Range("I132").Value = 1996.5
$ws.Range("K132").Value = 5989.5
$ws.Range("M132").Value = -3459.5

$ws.Range("H134").Value = 3088.9167
$ws.Range("I134").Value = 1312.75
$ws.Range("K134").Value = 3938.25
$ws.Range("M134").Value = -1403.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 35
$ws.Range("I38").Value = 36.25
$ws.Range("J38").Value = 30
$ws.Range("K38").Value = 108.75
$ws.Range("L38").Value = 90
$ws.Range("M38").Value = 238.25
$ws.Range("N38").Value = -784

$ws.Range("H97").Value = 1039.6666
$ws.Range("J97").Value = 542.6
$ws.Range("L97").Value = 1627.8
$ws.Range("N97").Value = -2619.8

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 13776.556
$ws.Range("J46").Value = 18333.334
$ws.Range("L46").Value = 18333.334
$ws.Range("N46").Value = -18645.334

$ws.Range("H80").Value = 3565.6155
$ws.Range("J80").Value = 3850.4
$ws.Range("L80").Value = 3850.4
$ws.Range("N80").Value = -5846.4

$ws.Range("H83").Value = 3565.6155
$ws.Range("J83").Value = 3850.4
$ws.Range("L83").Value = 19252
$ws.Range("N83").Value = -29236

$ws.Range("H122").Value = 12508550
$ws.Range("I122").Value = 31264124
$ws.Range("J122").Value = 4833.1665
$ws.Range("K122").Value = 93792372
$ws.Range("L122").Value = 14499.4995
$ws.Range("M122").Value = -93789922
$ws.Range("N122").Value = -19399.4995

$ws.Range("H126").Value = 4299.6
$ws.Range("I126").Value = 3249
$ws.Range("K126").Value = 9747
$ws.Range("M126").Value = -7277

$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 916.8
$ws.Range("I22").Value = 223.66667
$ws.Range("J22").Value = 1213.8572
$ws.Range("K22").Value = 223.66667
$ws.Range("L22").Value = 1213.8572
$ws.Range("M22").Value = 71.33332999999999
$ws.Range("N22").Value = -1803.8572

$ws.Range("H27").Value = 916.8
$ws.Range("I27").Value = 223.66667
$ws.Range("J27").Value = 1213.8572
$ws.Range("K27").Value = 223.66667
$ws.Range("L27").Value = 1213.8572
$ws.Range("M27").Value = -116.66667
$ws.Range("N27").Value = -1427.8572

$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 10000
$ws.Range("K39").Value = 10000
$ws.Range("M39").Value = -9540

$ws.Range("H40").Value = 4387.6665
$ws.Range("I40").Value = 3711.125
$ws.Range("J40").Value = 9800
$ws.Range("K40").Value = 3711.125
$ws.Range("L40").Value = 9800
$ws.Range("M40").Value = -3575.125
$ws.Range("N40").Value = -10072

$ws.Range("H55").Value = 725.2727
$ws.Range("I55").Value = 722.25
$ws.Range("J55").Value = 733.3333
$ws.Range("K55").Value = 722.25
$ws.Range("L55").Value = 733.3333
$ws.Range("M55").Value = -549.25
$ws.Range("N55").Value = -1079.3333

$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50450

$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51560

$ws.Range("H122").Value = 3716.923
$ws.Range("I122").Value = 3432
$ws.Range("K122").Value = 10296
$ws.Range("M122").Value = -7846

$ws.Range("H132").Value = 7336.8887
$ws.Range("I132").Value = 7330.304
$ws.Range("K132").Value = 21990.912
$ws.Range("M132").Value = -19460.912

$ws.Range("H136").Value = 2709.625
$ws.Range("I136").Value = 2589.9333
$ws.Range("J136").Value = 4505
$ws.Range("K136").Value = 7769.7999
$ws.Range("L136").Value = 13515
$ws.Range("M136").Value = -5219.7999
$ws.Range("N136").Value = -18615

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 333339940
$ws.Range("J11").Value = 9900
$ws.Range("L11").Value = 9900
$ws.Range("N11").Value = -10184

$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

